# Atualizando o arquivo XLSX
# A new match row was inserted as row 3 (COLOMBIA - PRIMERA A / Santa Fe vs
# Once Caldas). The rows that used to be row 3 (MEXICO - LIGA MX) and row 4
# (URUGUAY - PRIMERA DIVISION) are pushed down to rows 4 and 5 respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 3, shifting rows 3-4 down to 4-5.
$ws.Rows.Item(3).Insert()

# Populate the new row 3 with the inserted match data.
$ws.Range("A3").Value = "xnqclRX1"

# Date/Time columns are stored as plain text in this sheet, so force text
# formatting before assigning them to avoid Excel auto-converting them into
# date/time serial values, then drop back to the default style so no
# explicit number format is left behind on the cell.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "08/11/2024"
$ws.Range("B3").Style = "Normal"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "22:20"
$ws.Range("C3").Style = "Normal"

$ws.Range("D3").Value = "COLOMBIA - PRIMERA A"
$ws.Range("E3").Value = "Santa Fe"
$ws.Range("F3").Value = "Once Caldas"

$ws.Range("G3").Value = 2.1
$ws.Range("H3").Value = 3.1
$ws.Range("I3").Value = 3.8
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 1.91
$ws.Range("L3").Value = 4.75
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
$ws.Range("O3").Value = 1.53
$ws.Range("P3").Value = 2.38
$ws.Range("Q3").Value = 2.7
$ws.Range("R3").Value = 1.44
$ws.Range("S3").Value = 1.62
$ws.Range("T3").Value = 2.2
$ws.Range("U3").Value = 2.25
$ws.Range("V3").Value = 1.57
$ws.Range("W3").Value = 5.5
$ws.Range("X3").Value = 8.5
$ws.Range("Y3").Value = 10
$ws.Range("Z3").Value = 19
$ws.Range("AA3").Value = 21
$ws.Range("AB3").Value = 41
$ws.Range("AC3").Value = 6
$ws.Range("AD3").Value = 6
$ws.Range("AE3").Value = 21
$ws.Range("AF3").Value = 81
$ws.Range("AG3").Value = 201
$ws.Range("AH3").Value = 8
$ws.Range("AI3").Value = 17
$ws.Range("AJ3").Value = 15
$ws.Range("AK3").Value = 41
$ws.Range("AL3").Value = 41
$ws.Range("AM3").Value = 51
$ws.Range("AN3").Value = 4
$ws.Range("AO3").Value = 13
$ws.Range("AP3").Value = 29
$ws.Range("AQ3").Value = 41
$ws.Range("AR3").Value = 81
$ws.Range("AS3").Value = 301
$ws.Range("AT3").Value = 2.2
$ws.Range("AU3").Value = 9.5
$ws.Range("AV3").Value = 81
$ws.Range("AW3").Value = 5.5
$ws.Range("AX3").Value = 23
$ws.Range("AY3").Value = 41
$ws.Range("AZ3").Value = 81
$ws.Range("BA3").Value = 151
$ws.Range("BB3").Value = 351
$ws.Range("BC3").Value = 126
$ws.Range("BD3").Value = 126

Write-Host "Inserted new row 3 (xnqclRX1 / COLOMBIA - PRIMERA A) and shifted subsequent rows down."
